# Generate Report for Handoff
# A new localization entry (11ab0ace-ae72-45bd-875c-917c921c1f03) is being
# reported. It becomes the new "row 2" in every sheet, and the entry that used
# to live in row 2 (51c1ed14-94c1-4234-a34c-84ccff7282fe) is pushed down to a
# newly added row 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Remove the hyperlink currently anchored on A2 (it will be rebuilt below,
# once for the new row2 and once for the row3 that inherits the old data).
$existing = @()
foreach ($hl in $ws.Hyperlinks) { $existing += $hl }
foreach ($hl in $existing) { $hl.Delete() }

# Row 3 <- old row 2 payload (51c1ed14...)
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-30-21 06:30:36"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/51c1ed14-94c1-4234-a34c-84ccff7282fe.md", "", "", "51c1ed14-94c1-4234-a34c-84ccff7282fe.md") | Out-Null

# Row 2 <- new payload (11ab0ace...)
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-31-21 06:31:02"
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/11ab0ace-ae72-45bd-875c-917c921c1f03.md", "", "", "11ab0ace-ae72-45bd-875c-917c921c1f03.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$existing = @()
foreach ($hl in $ws.Hyperlinks) { $existing += $hl }
foreach ($hl in $existing) { $hl.Delete() }

# Row 3 <- old row 2 payload (51c1ed14...)
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-21 06:30:33"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/51c1ed14-94c1-4234-a34c-84ccff7282fe.md", "", "", "51c1ed14-94c1-4234-a34c-84ccff7282fe.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/51c1ed14-94c1-4234-a34c-84ccff7282fe.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6f3cd455a4265cf32a6075698c3a17d94f6f790e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/51c1ed14-94c1-4234-a34c-84ccff7282fe.cd45b03d24c9259dce136154e9fae89f337ebee1.zh-cn.xlf", "", "", "51c1ed14-94c1-4234-a34c-84ccff7282fe.cd45b03d24c9259dce136154e9fae89f337ebee1.zh-cn.xlf") | Out-Null

# Row 2 <- new payload (11ab0ace...)
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-21 06:30:58"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/11ab0ace-ae72-45bd-875c-917c921c1f03.md", "", "", "11ab0ace-ae72-45bd-875c-917c921c1f03.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/11ab0ace-ae72-45bd-875c-917c921c1f03.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e44847dc45ffc23a359455be2c1cd57585e9054b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/11ab0ace-ae72-45bd-875c-917c921c1f03.e44847dc45ffc23a359455be2c1cd57585e9054b.zh-cn.xlf", "", "", "11ab0ace-ae72-45bd-875c-917c921c1f03.e44847dc45ffc23a359455be2c1cd57585e9054b.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$existing = @()
foreach ($hl in $ws.Hyperlinks) { $existing += $hl }
foreach ($hl in $existing) { $hl.Delete() }

# Row 3 <- old row 2 payload (51c1ed14...)
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-21 06:30:36"
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/51c1ed14-94c1-4234-a34c-84ccff7282fe.md", "", "", "51c1ed14-94c1-4234-a34c-84ccff7282fe.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/51c1ed14-94c1-4234-a34c-84ccff7282fe.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc63aa1347a4c151d9dc753375871298d04a1aea/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/51c1ed14-94c1-4234-a34c-84ccff7282fe.cd45b03d24c9259dce136154e9fae89f337ebee1.de-de.xlf", "", "", "51c1ed14-94c1-4234-a34c-84ccff7282fe.cd45b03d24c9259dce136154e9fae89f337ebee1.de-de.xlf") | Out-Null

# Row 2 <- new payload (11ab0ace...)
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("E2").Value = "2016-03-21 06:31:02"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/11ab0ace-ae72-45bd-875c-917c921c1f03.md", "", "", "11ab0ace-ae72-45bd-875c-917c921c1f03.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/b020dbd473ff9902b0b9fd4eb70eb3b6c761a1f3/e2e/11ab0ace-ae72-45bd-875c-917c921c1f03.md", "", "", ".md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e44847dc45ffc23a359455be2c1cd57585e9054b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/11ab0ace-ae72-45bd-875c-917c921c1f03.e44847dc45ffc23a359455be2c1cd57585e9054b.de-de.xlf", "", "", "11ab0ace-ae72-45bd-875c-917c921c1f03.e44847dc45ffc23a359455be2c1cd57585e9054b.de-de.xlf") | Out-Null
